# This script applies the "456a3b4" data refresh to the 江西-漫展信息 workbook.
# For sheet "展览" (Exhibition) and sheet "全部类型" (All types):
#   - A new event row ("宜春·ML宅舞奖金赛-宜春万达赛区") is inserted right
#     after the "九江·动漫畅想（取消）" row (i.e. becomes the new row 12),
#     pushing every subsequent row down by one.
#   - Several "想去人数" (interested-count) values are refreshed (bumped up).
# For sheet "演出" (Performance) only the interested-count of the CrossingX
# row is refreshed.
#
# NOTE: this engine's PowerShell does not bind *named* function parameters
# correctly, so all helper functions below use purely positional parameters.

$wb = $excel.ActiveWorkbook

function Insert-NewRow12 {
    param($ws)

    # Insert a new blank row at position 12 (shifts rows 12.. down by one)
    $ws.Rows.Item(12).Insert()

    # Copy the formatting from row 11 (the row just above) onto the new row 12
    # so the new row looks consistent with the rest of the table.
    $ws.Range("A11:I11").Copy()
    $ws.Range("A12:I12").PasteSpecial(-4122)   # xlPasteFormats

    # Populate the new row's contents.
    $ws.Range("A12").Value = 11

    $ws.Range("B12").NumberFormat = "@"
    $ws.Range("B12").Value = "2024-08-11"
    $ws.Range("B12").Style = "Normal"

    $ws.Range("C12").Value = "宜春·ML宅舞奖金赛-宜春万达赛区"
    $ws.Range("D12").Value = "宜阳东大道6号 宜春万达广场"

    $ws.Range("E12").NumberFormat = "@"
    $ws.Range("E12").Value = "2024.08.11 14:00-08.11 19:00"
    $ws.Range("E12").Style = "Normal"

    $ws.Range("F12").Value = 3
    $ws.Range("G12").Value = 30
    $ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=90446"
    $ws.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202408/f6Rm6dm61722996358480.jpeg"
}

function Fix-SequenceColumn {
    param($ws, $lastRow)
    # Column A simply holds (row number - 1); after the insert the shifted
    # rows still carry their old numbers, so renumber rows 12..lastRow.
    for ($r = 12; $r -le $lastRow; $r++) {
        $ws.Range("A$r").Value = $r - 1
    }
}

function Set-F {
    param($ws, $row, $value)
    $ws.Range("F$row").Value = $value
}

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

Insert-NewRow12 $ws1
Fix-SequenceColumn $ws1 26

# Refresh "want to go" counts on the rows that kept their position.
Set-F $ws1 2  1483
Set-F $ws1 3  145
Set-F $ws1 4  1771
Set-F $ws1 6  148
Set-F $ws1 7  660
Set-F $ws1 8  36
Set-F $ws1 10 555

# Refresh "want to go" counts on the rows that shifted down by one.
Set-F $ws1 13 81     # 抚州·逆光ZERO动漫游戏展
Set-F $ws1 14 152    # 萍乡·夏花Flower·2024夏季国漫展
Set-F $ws1 15 22     # 新余·逆光ZERO动漫游戏展 (unchanged)
Set-F $ws1 16 122    # 上饶·次元重现夏日嘉年华（取消）
Set-F $ws1 17 72     # 乐平·CY境界次元第三届动漫游戏庆典
Set-F $ws1 18 105    # 南昌·CM03·配音演员孙路路专场见面会 (unchanged)
Set-F $ws1 19 4894   # 南昌·CM03动漫游戏博览会
Set-F $ws1 20 48     # 九江·如梦令国潮动漫节 (unchanged)
Set-F $ws1 21 828    # 南昌·第四届龙年动漫展——暑假最后的狂欢
Set-F $ws1 22 111    # 赣州·第五人格only
Set-F $ws1 23 2227   # 南昌·Sunflower Garden动漫游戏展
Set-F $ws1 24 70     # 南昌·第一届哥布林动漫游戏展——开学季&贺中秋 (unchanged)
Set-F $ws1 25 19     # 南昌·Aud中秋动漫嘉年华
Set-F $ws1 26 2081   # 南昌·萌卡动漫展

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performance)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
Set-F $ws2 2 78      # 南昌·CrossingX意次元｜乐队番ONLY

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

Insert-NewRow12 $ws4
Fix-SequenceColumn $ws4 28

Set-F $ws4 2  1483
Set-F $ws4 3  145
Set-F $ws4 4  1771
Set-F $ws4 6  148
Set-F $ws4 7  660
Set-F $ws4 8  36
Set-F $ws4 10 555

Set-F $ws4 13 81     # 抚州·逆光ZERO动漫游戏展
Set-F $ws4 14 152    # 萍乡·夏花Flower·2024夏季国漫展
Set-F $ws4 15 22     # 新余·逆光ZERO动漫游戏展 (unchanged)
Set-F $ws4 16 122    # 上饶·次元重现夏日嘉年华（取消）
Set-F $ws4 17 72     # 乐平·CY境界次元第三届动漫游戏庆典
Set-F $ws4 18 105    # 南昌·CM03·配音演员孙路路专场见面会 (unchanged)
Set-F $ws4 19 4894   # 南昌·CM03动漫游戏博览会
Set-F $ws4 20 78     # 南昌·CrossingX意次元｜乐队番ONLY
Set-F $ws4 21 48     # 九江·如梦令国潮动漫节 (unchanged)
Set-F $ws4 22 2      # 南昌·【8月24日】滑稽互动狂欢大作战《欢乐小丑嘉年华》 (unchanged)
Set-F $ws4 23 828    # 南昌·第四届龙年动漫展——暑假最后的狂欢
Set-F $ws4 24 111    # 赣州·第五人格only
Set-F $ws4 25 2227   # 南昌·Sunflower Garden动漫游戏展
Set-F $ws4 26 70     # 南昌·第一届哥布林动漫游戏展——开学季&贺中秋 (unchanged)
Set-F $ws4 27 19     # 南昌·Aud中秋动漫嘉年华
Set-F $ws4 28 2081   # 南昌·萌卡动漫展

Write-Host "Edit complete."
